$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 345-346, pushing the existing rows 345..376 down to 347..378
$ws.Range("A345:A346").EntireRow.Insert()

# New "Primera" record (row 345) - most recent week added at the top of this block
$ws.Range("A345").Value = 1
$ws.Range("B345").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C345").Value = "Arica y Parinacota"
$ws.Range("D345").Value = 44858
$ws.Range("E345").Value = 15
$ws.Range("F345").Value = 100112043
$ws.Range("G345").Value = "Pepino ensalada"
$ws.Range("H345").Value = "Sin especificar"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 250
$ws.Range("K345").Value = 13000
$ws.Range("L345").Value = 14000
$ws.Range("M345").Value = 13600
$ws.Range("N345").Value = "$/caja 70 unidades"
$ws.Range("O345").Value = "Región de Arica y Parinacota"
$ws.Range("P345").Value = 194
$ws.Range("Q345").Value = 70
$ws.Range("R345").Value = "Hortaliza"

# New "Segunda" record (row 346)
$ws.Range("A346").Value = 1
$ws.Range("B346").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C346").Value = "Arica y Parinacota"
$ws.Range("D346").Value = 44858
$ws.Range("E346").Value = 15
$ws.Range("F346").Value = 100112043
$ws.Range("G346").Value = "Pepino ensalada"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Segunda"
$ws.Range("J346").Value = 350
$ws.Range("K346").Value = 9000
$ws.Range("L346").Value = 10000
$ws.Range("M346").Value = 9429
$ws.Range("N346").Value = "$/caja 100 unidades"
$ws.Range("O346").Value = "Región de Arica y Parinacota"
$ws.Range("P346").Value = 94
$ws.Range("Q346").Value = 100
$ws.Range("R346").Value = "Hortaliza"
